# "Update diagrammes classes bis"
# Target shape: slide 1, "Rectangle 55" (shape id 56) - the (empty) methods
# compartment of the "Messagerie" class box. We grow its height and fill in
# the method signatures (mirroring the "Serveur" class's method box).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 55")

# --- grow the box to fit the new text -------------------------------------
# cy: 599369 EMU -> 754944 EMU (off/cx stay the same); COM Height/Width are
# in points (1 pt = 12700 EMU) and are single-precision, so nudge slightly
# inside the rounding band that lands on exactly 754944 EMU.
$shp.Height = 59.44443

# --- fill in the method list ------------------------------------------------
$tr = $shp.TextFrame.TextRange

function Format-Run($rng) {
    $rng.LanguageID = "en-US"
    $rng.Font.Size = 9
    $rng.Font.Color.RGB = 0
}

# Paragraph 1: "+ connect(String id)"
$run = $tr.InsertBefore("+ connect(String id)")
Format-Run $run

# Paragraph 2: "+ send(String message, String id)"
$run = $run.InsertAfter("`r+ send(String message, String id)")
Format-Run $run

# Paragraph 3: "+ bye (String id)"
$run = $run.InsertAfter("`r+ bye (String id)")
Format-Run $run

# Paragraph 4: "+ who() : " & "Liste" & "<String>"
$run = $run.InsertAfter("`r+ who() : ")
Format-Run $run
$run = $run.InsertAfter("Liste")
Format-Run $run
$run = $run.InsertAfter("<String>")
Format-Run $run

# Paragraph 5: "+ " & "getMessages" & "(" & "int" & " " & "nbMessage" & ")"
$run = $run.InsertAfter("`r+ ")
Format-Run $run
$run = $run.InsertAfter("getMessages")
Format-Run $run
$run = $run.InsertAfter("(")
Format-Run $run
$run = $run.InsertAfter("int")
Format-Run $run
$run = $run.InsertAfter(" ")
Format-Run $run
$run = $run.InsertAfter("nbMessage")
Format-Run $run
$run = $run.InsertAfter(")")
Format-Run $run
